$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Calr"
$ws.Range("C2").Value = "Itga3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 99.883077
$ws.Range("H2").Value = 299.649231
$ws.Range("I2").Value = 0.3917580942718377
$ws.Range("J2").Value = 0.3917580942718377
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.269168666666666
$ws.Range("N2").Value = 24.807506
$ws.Range("O2").Value = 0.671680253471746
$ws.Range("P2").Value = 0.671680253471746
$ws.Range("Q2").Value = 825.9500106586539
$ws.Range("R2").Value = 7433.550095927884
$ws.Range("S2").Value = 0.2631361760601161
$ws.Range("T2").Value = 0.2631361760601161

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Calr"
$ws.Range("C3").Value = "Itga3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 99.883077
$ws.Range("H3").Value = 299.649231
$ws.Range("I3").Value = 0.3917580942718377
$ws.Range("J3").Value = 0.3917580942718377
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.121294
$ws.Range("N3").Value = 0.363882
$ws.Range("O3").Value = 0.009852354928133683
$ws.Range("P3").Value = 0.009852354928133683
$ws.Range("Q3").Value = 12.115217941638
$ws.Range("R3").Value = 109.036961474742
$ws.Range("S3").Value = 0.0038597397907354
$ws.Range("T3").Value = 0.0038597397907354

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Calr"
$ws.Range("C4").Value = "Itga3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 99.883077
$ws.Range("H4").Value = 299.649231
$ws.Range("I4").Value = 0.3917580942718377
$ws.Range("J4").Value = 0.3917580942718377
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.920705666666667
$ws.Range("N4").Value = 11.762117
$ws.Range("O4").Value = 0.3184673916001203
$ws.Range("P4").Value = 0.3184673916001203
$ws.Range("Q4").Value = 391.612145998003
$ws.Range("R4").Value = 3524.509313982027
$ws.Range("S4").Value = 0.1247621784209862
$ws.Range("T4").Value = 0.1247621784209862

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Calr"
$ws.Range("C5").Value = "Itga3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 124.0161413333333
$ws.Range("H5").Value = 372.048424
$ws.Range("I5").Value = 0.4864119993789693
$ws.Range("J5").Value = 0.4864119993789694
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.269168666666666
$ws.Range("N5").Value = 24.807506
$ws.Range("O5").Value = 0.671680253471746
$ws.Range("P5").Value = 0.671680253471746
$ws.Range("Q5").Value = 1025.510390074505
$ws.Range("R5").Value = 9229.593510670544
$ws.Range("S5").Value = 0.3267133350345648
$ws.Range("T5").Value = 0.3267133350345649

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Calr"
$ws.Range("C6").Value = "Itga3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 124.0161413333333
$ws.Range("H6").Value = 372.048424
$ws.Range("I6").Value = 0.4864119993789693
$ws.Range("J6").Value = 0.4864119993789694
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.121294
$ws.Range("N6").Value = 0.363882
$ws.Range("O6").Value = 0.009852354928133683
$ws.Range("P6").Value = 0.009852354928133683
$ws.Range("Q6").Value = 15.04241384688533
$ws.Range("R6").Value = 135.381724621968
$ws.Range("S6").Value = 0.004792303659184746
$ws.Range("T6").Value = 0.004792303659184747

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Calr"
$ws.Range("C7").Value = "Itga3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 124.0161413333333
$ws.Range("H7").Value = 372.048424
$ws.Range("I7").Value = 0.4864119993789693
$ws.Range("J7").Value = 0.4864119993789694
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.920705666666667
$ws.Range("N7").Value = 11.762117
$ws.Range("O7").Value = 0.3184673916001203
$ws.Range("P7").Value = 0.3184673916001203
$ws.Range("Q7").Value = 486.2307880837342
$ws.Range("R7").Value = 4376.077092753608
$ws.Range("S7").Value = 0.1549063606852197
$ws.Range("T7").Value = 0.1549063606852197

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Calr"
$ws.Range("C8").Value = "Itga3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 31.06188766666667
$ws.Range("H8").Value = 93.185663
$ws.Range("I8").Value = 0.121829906349193
$ws.Range("J8").Value = 0.121829906349193
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.269168666666666
$ws.Range("N8").Value = 24.807506
$ws.Range("O8").Value = 0.671680253471746
$ws.Range("P8").Value = 0.671680253471746
$ws.Range("Q8").Value = 256.8559882207197
$ws.Range("R8").Value = 2311.703893986478
$ws.Range("S8").Value = 0.08183074237706502
$ws.Range("T8").Value = 0.08183074237706503

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Calr"
$ws.Range("C9").Value = "Itga3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 31.06188766666667
$ws.Range("H9").Value = 93.185663
$ws.Range("I9").Value = 0.121829906349193
$ws.Range("J9").Value = 0.121829906349193
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.121294
$ws.Range("N9").Value = 0.363882
$ws.Range("O9").Value = 0.009852354928133683
$ws.Range("P9").Value = 0.009852354928133683
$ws.Range("Q9").Value = 3.767620602640667
$ws.Range("R9").Value = 33.908585423766
$ws.Range("S9").Value = 0.001200311478213537
$ws.Range("T9").Value = 0.001200311478213537

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Calr"
$ws.Range("C10").Value = "Itga3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 31.06188766666667
$ws.Range("H10").Value = 93.185663
$ws.Range("I10").Value = 0.121829906349193
$ws.Range("J10").Value = 0.121829906349193
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.920705666666667
$ws.Range("N10").Value = 11.762117
$ws.Range("O10").Value = 0.3184673916001203
$ws.Range("P10").Value = 0.3184673916001203
$ws.Range("Q10").Value = 121.7845189920634
$ws.Range("R10").Value = 1096.060670928571
$ws.Range("S10").Value = 0.03879885249391442
$ws.Range("T10").Value = 0.03879885249391443

